$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.994.33'
$ws.Range("D3").Value = '2.171.96'
$ws.Range("E3").Value = '  -2.34%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.35'
$ws.Range("E5").Value = '  -3.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.612'
$ws.Range("E6").Value = '  -2.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.63'
$ws.Range("E7").Value = '  -5.61%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.572'
$ws.Range("E9").Value = '  -6.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.88'
$ws.Range("E10").Value = '  -8.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0928'
$ws.Range("E11").Value = '  -3.63%  '
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.74'
$ws.Range("E13").Value = '  -5.20%  '
$ws.Range("D14").Value = '2.495.43'
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.82'
$ws.Range("E15").Value = '  -3.08%  '
$ws.Range("E16").Value = '  -4.21%  '
$ws.Range("D17").Value = '2.166.81'
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("D18").Value = '40.866.55'
$ws.Range("E18").Value = '  -2.64%  '
$ws.Range("E19").Value = '  -7.32%  '
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("E21").Value = '  -4.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.81'
$ws.Range("E22").Value = '  -4.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '225.05'
$ws.Range("E23").Value = '  -2.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").Value = '  -7.69%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("E26").Value = '  -5.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.54'
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  -3.99%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '165.82'
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.82'
$ws.Range("E31").Value = '  -3.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '30.88'
$ws.Range("E32").Value = '  +4.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0773'
$ws.Range("E33").Value = '  -3.19%  '
$ws.Range("E34").Value = '  -8.03%  '
$ws.Range("E35").Value = '  -3.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.104'
$ws.Range("E36").Value = '  -9.25%  '
$ws.Range("E37").Value = '  -3.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0287'
$ws.Range("E38").Value = '  -5.00%  '
$ws.Range("E39").Value = '  -5.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.06'
$ws.Range("E40").Value = '  -3.93%  '
$ws.Range("E41").Value = '  -4.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '59.74'
$ws.Range("E42").Value = '  -7.29%  '
$ws.Range("E43").Value = '  -5.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.28'
$ws.Range("E44").Value = '  -5.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0972'
$ws.Range("E45").Value = '  -3.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '98.90'
$ws.Range("E46").Value = '  -5.61%  '
$ws.Range("E47").Value = '  -3.45%  '
$ws.Range("E48").Value = '  -3.16%  '
$ws.Range("E49").Value = '  -7.81%  '
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("D51").Value = '2.373.08'
$ws.Range("E51").Value = '  -2.42%  '
